$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.256.54'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '3.471.67'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '593.25'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '178.37'
$ws.Range("E6").Value = '  +4.02%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.473.74'
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '0.137'
$ws.Range("E10").Value = '  +5.64%  '
$ws.Range("E11").Value = '  -2.46%  '
$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '4.072.83'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("D14").Value = '31.95'
$ws.Range("E14").Value = '  +11.36%  '
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").Value = '67.290.82'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '3.466.15'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").Value = '14.26'
$ws.Range("E20").Value = '  +1.86%  '
$ws.Range("D21").Value = '387.78'
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("E22").Value = '  -0.80%  '
$ws.Range("D23").Value = '72.80'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '5.72'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = '0.175'
$ws.Range("E29").Value = '  -3.07%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = '6.16'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '23.47'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("D35").Value = '7.37'
$ws.Range("E35").Value = '  +0.97%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("D38").Value = '163.79'
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").Value = '2.72'
$ws.Range("E41").Value = '  +6.51%  '
$ws.Range("D42").Value = '6.85'
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '4.60'
$ws.Range("E43").Value = '  -1.07%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.823.46'
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").Value = '26.10'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").Value = '0.0721'
$ws.Range("E46").Value = '  -2.29%  '
$ws.Range("D47").Value = '26.48'
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("D48").Value = '41.46'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("D49").Value = '0.0297'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").Value = '336.05'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '1.04'
$ws.Range("E51").Value = '  -2.58%  '
